# Redo naming of mixres units
#
# The old "area_mixres" sheet (the original, now-superseded stats) is
# dropped entirely, and the sheet that was "area_mixres_new" takes over
# as the (renamed) "area_mixre" sheet. All other sheets (area_hires,
# area_lores, area_pop_sum) simply shift left by one position, keeping
# their own names/content unchanged.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$oldSheet = $wb.Worksheets.Item("area_mixres")
$oldSheet.Delete() | Out-Null

$newSheet = $wb.Worksheets.Item("area_mixres_new")
$newSheet.Name = "area_mixre"
